$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update mac_address column (C) for existing rows 2-21 to new dash-format MAC values
$ws.Range("C2").Value = "8C-16-45-5A-5D-0D"
$ws.Range("C3").Value = "8C-16-45-88-E1-0D"
$ws.Range("C4").Value = "00-FF-D3-E3-9A-27"
$ws.Range("C5").Value = "8C-16-45-5A-62-41"
$ws.Range("C6").Value = "E8-6A-64-1D-75-E4"
$ws.Range("C7").Value = "8C-16-45-FA-94-B7"
$ws.Range("C8").Value = "8C-16-45-1A-0F-62"
$ws.Range("C9").Value = "E8-6A-64-1C-52-6E"
$ws.Range("C10").Value = "48-51-B7-10-35-A6"
$ws.Range("C11").Value = "8C-16-45-38-F3-F3"
$ws.Range("C12").Value = "D4-3D-7E-58-CC-45"
$ws.Range("C13").Value = "8C-16-45-5A-5D-96"
$ws.Range("C14").Value = "8C-16-45-5A-5D-8E"
$ws.Range("C15").Value = "8C-16-45-33-A5-5F"
$ws.Range("C16").Value = "3C-95-09-F9-EA-DF"
$ws.Range("C17").Value = "8C-16-45-88-E7-0B"
$ws.Range("C18").Value = "B4-69-21-5A-DB-C4"
$ws.Range("C19").Value = "E8-6A-64-1D-48-B7"
$ws.Range("C20").Value = "8C-16-45-59-69-09 "
$ws.Range("C21").Value = "98-E7-F4-30-16-5A "

# Add new rows 22-30 for Machine 21-29
$ws.Range("A22").Value = 10021
$ws.Range("B22").Value = "Machine 21"
$ws.Range("C22").Value = "38-BA-F8-53-C7-8F"
$ws.Range("D22").Value = "FB5962911653"
$ws.Range("E22").Value = "192.168.0.874"
$ws.Range("F22").Value = 1001
$ws.Range("G22").Value = "eng"
$ws.Range("H22").Value = $True
$ws.Range("I22").Value = "superadmin"
$ws.Range("J22").Value = "now()"

$ws.Range("A23").Value = 10022
$ws.Range("B23").Value = "Machine 22"
$ws.Range("C23").Value = "E8-6A-64-1C-58-C2"
$ws.Range("D23").Value = "FB5962911654"
$ws.Range("E23").Value = "192.168.0.721"
$ws.Range("F23").Value = 1001
$ws.Range("G23").Value = "eng"
$ws.Range("H23").Value = $True
$ws.Range("I23").Value = "superadmin"
$ws.Range("J23").Value = "now()"

$ws.Range("A24").Value = 10023
$ws.Range("B24").Value = "Machine 23"
$ws.Range("C24").Value = "E4-A4-71-CE-BA-93"
$ws.Range("D24").Value = "FB5962911655"
$ws.Range("E24").Value = "192.168.0.841"
$ws.Range("F24").Value = 1001
$ws.Range("G24").Value = "eng"
$ws.Range("H24").Value = $True
$ws.Range("I24").Value = "superadmin"
$ws.Range("J24").Value = "now()"

$ws.Range("A25").Value = 10024
$ws.Range("B25").Value = "Machine 24"
$ws.Range("C25").Value = "54-E1-AD-EA-30-C9"
$ws.Range("D25").Value = "FB5962911656"
$ws.Range("E25").Value = "192.168.0.186"
$ws.Range("F25").Value = 1001
$ws.Range("G25").Value = "eng"
$ws.Range("H25").Value = $True
$ws.Range("I25").Value = "superadmin"
$ws.Range("J25").Value = "now()"

$ws.Range("A26").Value = 10025
$ws.Range("B26").Value = "Machine 25"
$ws.Range("C26").Value = "8C-16-45-65-DD-40"
$ws.Range("D26").Value = "FB5962911657"
$ws.Range("E26").Value = "192.168.0.627"
$ws.Range("F26").Value = 1001
$ws.Range("G26").Value = "eng"
$ws.Range("H26").Value = $True
$ws.Range("I26").Value = "superadmin"
$ws.Range("J26").Value = "now()"

$ws.Range("A27").Value = 10026
$ws.Range("B27").Value = "Machine 26"
$ws.Range("C27").Value = "58-20-B1-D6-C3-BE"
$ws.Range("D27").Value = "FB5962911658"
$ws.Range("E27").Value = "192.168.0.879"
$ws.Range("F27").Value = 1001
$ws.Range("G27").Value = "eng"
$ws.Range("H27").Value = $True
$ws.Range("I27").Value = "superadmin"
$ws.Range("J27").Value = "now()"

$ws.Range("A28").Value = 10027
$ws.Range("B28").Value = "Machine 27"
$ws.Range("C28").Value = "8C-16-45-38-F0-25"
$ws.Range("D28").Value = "FB5962911659"
$ws.Range("E28").Value = "192.168.0.628"
$ws.Range("F28").Value = 1001
$ws.Range("G28").Value = "eng"
$ws.Range("H28").Value = $True
$ws.Range("I28").Value = "superadmin"
$ws.Range("J28").Value = "now()"

$ws.Range("A29").Value = 10028
$ws.Range("B29").Value = "Machine 28"
$ws.Range("C29").Value = "6C-88-14-AC-EF-55"
$ws.Range("D29").Value = "FB5962911661"
$ws.Range("E29").Value = "192.168.0.306"
$ws.Range("F29").Value = 1001
$ws.Range("G29").Value = "eng"
$ws.Range("H29").Value = $True
$ws.Range("I29").Value = "superadmin"
$ws.Range("J29").Value = "now()"

$ws.Range("A30").Value = 10029
$ws.Range("B30").Value = "Machine 29"
$ws.Range("C30").Value = "3C-6A-A7-C0-DF-27"
$ws.Range("D30").Value = "FB5962911662"
$ws.Range("E30").Value = "192.168.0.355"
$ws.Range("F30").Value = 1001
$ws.Range("G30").Value = "eng"
$ws.Range("H30").Value = $True
$ws.Range("I30").Value = "superadmin"
$ws.Range("J30").Value = "now()"

# Adjust column C width to fit the new (longer) mac_address values
$ws.Columns.Item(3).ColumnWidth = 16.1666666666666668

# Move the selection to the next empty row, as left by the editor
$ws.Range("A31:XFD1048576").Select()